# Adapt column header formatting to respective input file names:
#   *_old -> *_FV2210, *_new -> *_FV2304
# and turn the header range into an Excel Table with a frozen header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row's shared-string labels -----------------------
$lastCol = 21   # A .. U
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $v = $cell.Value2
    if ($null -ne $v) {
        if ($v -like "*_old") {
            $cell.Value2 = ($v -replace "_old$", "_FV2210")
        } elseif ($v -like "*_new") {
            $cell.Value2 = ($v -replace "_new$", "_FV2304")
        }
    }
}

# --- 2. Freeze the header row ----------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Turn the used range into an Excel Table -----------------------------
$dim = $ws.Range("A1:U58")
$lo = $ws.ListObjects.Add(1, $dim, $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""
